$wb = $excel.ActiveWorkbook

# Values updated per the Sheets scheduled-runner refresh (market price recalculation).
# Each block targets one worksheet and writes the refreshed currentAveragePrice /
# LevePrice / LeveProfit columns (H, I, J, K, L, M, N) for the affected leve rows.

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3941.8
$ws.Range("I51").Value = 3400
$ws.Range("J51").Value = 4002
$ws.Range("K51").Value = 3400
$ws.Range("L51").Value = 4002
$ws.Range("M51").Value = -2916
$ws.Range("N51").Value = -4970
$ws.Range("H74").Value = 6625.0303
$ws.Range("I74").Value = 3948.625
$ws.Range("K74").Value = 3948.625
$ws.Range("M74").Value = -3012.625
$ws.Range("H77").Value = 6625.0303
$ws.Range("I77").Value = 3948.625
$ws.Range("K77").Value = 19743.125
$ws.Range("M77").Value = -15063.125
$ws.Range("H111").Value = 593.2222
$ws.Range("I111").Value = 588.375
$ws.Range("K111").Value = 1765.125
$ws.Range("M111").Value = 1301.875
$ws.Range("H125").Value = 2305.8
$ws.Range("I125").Value = 2337
$ws.Range("J125").Value = 2274.6
$ws.Range("K125").Value = 21033
$ws.Range("L125").Value = 20471.4
$ws.Range("M125").Value = -18573
$ws.Range("N125").Value = -25391.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3497.5
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = $null
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = $null
$ws.Range("N79").Value = 0
$ws.Range("H97").Value = 697.1429000000001
$ws.Range("I97").Value = 697.1429000000001
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 697.1429000000001
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = $null
$ws.Range("N97").Value = -201.1429000000001
$ws.Range("H110").Value = 444.06897
$ws.Range("I110").Value = 442.67856
$ws.Range("J110").Value = 483
$ws.Range("K110").Value = 442.67856
$ws.Range("L110").Value = 483
$ws.Range("M110").Value = 1602.32144
$ws.Range("N110").Value = -4573
$ws.Range("H132").Value = 3674.487
$ws.Range("I132").Value = 3718.0527
$ws.Range("K132").Value = 11154.1581
$ws.Range("M132").Value = -8624.158100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 53500
$ws.Range("J35").Value = 53500
$ws.Range("L35").Value = 53500
$ws.Range("N35").Value = -54120
$ws.Range("H88").Value = 50343
$ws.Range("J88").Value = 50343
$ws.Range("L88").Value = 50343
$ws.Range("N88").Value = -51155
$ws.Range("H91").Value = 50343
$ws.Range("J91").Value = 50343
$ws.Range("L91").Value = 50343
$ws.Range("N91").Value = -53151
$ws.Range("H105").Value = 2010.579
$ws.Range("I105").Value = 2093.9333
$ws.Range("J105").Value = 1698
$ws.Range("K105").Value = 2093.9333
$ws.Range("L105").Value = 1698
$ws.Range("M105").Value = -346.9333000000001
$ws.Range("N105").Value = -5192

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = $null
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = $null
$ws.Range("N67").Value = 0
$ws.Range("H70").Value = 74925.336
$ws.Range("J70").Value = 74925.336
$ws.Range("L70").Value = 74925.336
$ws.Range("N70").Value = -75555.336
$ws.Range("H73").Value = 74925.336
$ws.Range("J73").Value = 74925.336
$ws.Range("L73").Value = 74925.336
$ws.Range("N73").Value = -77109.336
$ws.Range("H86").Value = 4472
$ws.Range("I86").Value = 3962.6667
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 3962.6667
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -2839.6667
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 4472
$ws.Range("I89").Value = 3962.6667
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 19813.3335
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -14197.3335
$ws.Range("N89").Value = -41232
$ws.Range("H105").Value = 2196.6667
$ws.Range("I105").Value = 2495
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 2495
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = -748
$ws.Range("N105").Value = -5094
$ws.Range("H132").Value = 1615.1428
$ws.Range("I132").Value = 1362.8422
$ws.Range("K132").Value = 4088.5266
$ws.Range("M132").Value = -1558.5266
$ws.Range("H134").Value = 2757.2307
$ws.Range("I134").Value = 2570.3333
$ws.Range("K134").Value = 7710.999899999999
$ws.Range("M134").Value = -5175.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 6936.25
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 7784.2856
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 23352.8568
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -24974.8568
$ws.Range("H72").Value = 6936.25
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 7784.2856
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 70058.5704
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -78170.5704

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 45000
$ws.Range("J15").Value = 45000
$ws.Range("L15").Value = 45000
$ws.Range("N15").Value = -45576
$ws.Range("H80").Value = 8224
$ws.Range("J80").Value = 9949.5
$ws.Range("L80").Value = 9949.5
$ws.Range("N80").Value = -11945.5
$ws.Range("H81").Value = 45000
$ws.Range("J81").Value = 45000
$ws.Range("L81").Value = 45000
$ws.Range("N81").Value = -46996
$ws.Range("H83").Value = 8224
$ws.Range("J83").Value = 9949.5
$ws.Range("L83").Value = 49747.5
$ws.Range("N83").Value = -59731.5
$ws.Range("H84").Value = 45000
$ws.Range("J84").Value = 45000
$ws.Range("L84").Value = 135000
$ws.Range("N84").Value = -144984
$ws.Range("H102").Value = 1682.3784
$ws.Range("I102").Value = 1061.8
$ws.Range("K102").Value = 1061.8
$ws.Range("M102").Value = 560.2
$ws.Range("H113").Value = 2721.842
$ws.Range("I113").Value = 1081.3
$ws.Range("J113").Value = 4544.6665
$ws.Range("K113").Value = 1081.3
$ws.Range("L113").Value = 4544.6665
$ws.Range("M113").Value = 1088.7
$ws.Range("N113").Value = -8884.666499999999
$ws.Range("H126").Value = 2604.4814
$ws.Range("J126").Value = 2637.5
$ws.Range("L126").Value = 7912.5
$ws.Range("N126").Value = -12852.5
$ws.Range("H132").Value = 2201.2122
$ws.Range("I132").Value = 2234.913
$ws.Range("K132").Value = 6704.739
$ws.Range("M132").Value = -4174.739

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3569.7334
$ws.Range("I7").Value = 3340.5454
$ws.Range("K7").Value = 3340.5454
$ws.Range("M7").Value = -3228.5454
$ws.Range("H29").Value = 2000000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
$ws.Range("H43").Value = 602400
$ws.Range("J43").Value = 1000000
$ws.Range("L43").Value = 1000000
$ws.Range("N43").Value = -1000386
$ws.Range("H76").Value = 9644
$ws.Range("J76").Value = 9288
$ws.Range("L76").Value = 9288
$ws.Range("N76").Value = -9964
$ws.Range("H79").Value = 9644
$ws.Range("J79").Value = 9288
$ws.Range("L79").Value = 9288
$ws.Range("N79").Value = -11628
$ws.Range("H100").Value = 7361.25
$ws.Range("I100").Value = 2604.25
$ws.Range("J100").Value = 16875.25
$ws.Range("K100").Value = 2604.25
$ws.Range("L100").Value = 16875.25
$ws.Range("M100").Value = -2063.25
$ws.Range("N100").Value = -17957.25
$ws.Range("H126").Value = 3569.7334
$ws.Range("I126").Value = 3340.5454
$ws.Range("K126").Value = 10021.6362
$ws.Range("M126").Value = -7551.636200000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = $null
$ws.Range("N42").Value = 0
$ws.Range("H62").Value = 119125
$ws.Range("I62").Value = 900000
$ws.Range("J62").Value = 7571.4287
$ws.Range("K62").Value = 900000
$ws.Range("L62").Value = 7571.4287
$ws.Range("M62").Value = -899376
$ws.Range("N62").Value = -8819.4287
$ws.Range("H65").Value = 119125
$ws.Range("I65").Value = 900000
$ws.Range("J65").Value = 7571.4287
$ws.Range("K65").Value = 4500000
$ws.Range("L65").Value = 37857.14350000001
$ws.Range("M65").Value = -4496880
$ws.Range("N65").Value = -44097.14350000001
$ws.Range("H68").Value = 60180.332
$ws.Range("H71").Value = 60180.332
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = $null
$ws.Range("N82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = $null
$ws.Range("N85").Value = 0
$ws.Range("H86").Value = 49986.668
$ws.Range("J86").Value = 49986.668
$ws.Range("L86").Value = 49986.668
$ws.Range("N86").Value = -52232.668
$ws.Range("H89").Value = 49986.668
$ws.Range("J89").Value = 49986.668
$ws.Range("L89").Value = 249933.34
$ws.Range("N89").Value = -261165.34
$ws.Range("H100").Value = 1710.3
$ws.Range("I100").Value = 1756.1538
$ws.Range("J100").Value = 1625.1428
$ws.Range("K100").Value = 3512.3076
$ws.Range("L100").Value = 3250.2856
$ws.Range("M100").Value = -2971.3076
$ws.Range("N100").Value = -4332.2856
$ws.Range("H132").Value = 2700.5208
$ws.Range("I132").Value = 1860.4722
$ws.Range("K132").Value = 5581.4166
$ws.Range("M132").Value = -3051.4166
